# REQUIREMENTS.docx edit:
#  - finish the "make ga|meplay" sentence so it reads "...will make
#    gameplay more convenient..." and keep the closing sentence about
#    the rest of the document.
#  - append a new block of section headers/placeholders after the
#    INTRODUCTION paragraph: a "GIT TAG GOES HERE..." placeholder line,
#    a SYSTEM DIAGRAM header + caption, and ACTIONS AND THEIR SCENARIOS /
#    PLATFORM / SUMMARY headers. The trailing "_GoBack" bookmark (Word's
#    auto "last edit" marker) ends up right after the new SUMMARY text.

$d = $word.ActiveDocument

# The existing "_GoBack" bookmark currently sits between "...will make
# ga" and "meplay more convenient...".  We'll delete everything from
# that point on (the tail of the Introduction paragraph) and rebuild it,
# then keep growing the document by always inserting new text
# immediately *before* the (self-relocating) bookmark -- this is exactly
# how typing new content right at the bookmark's position behaves, so
# the bookmark naturally ends up past everything we add, landing at the
# very end of the new SUMMARY heading, matching the target layout.

$full = $d.Content.Text
$tailStart = $full.IndexOf("meplay")
$tailStart = $full.IndexOf("meplay", $tailStart + 1)
$tailEnd = $full.Length
$tailRange = $d.Range($tailStart, $tailEnd - 1)
$tailRange.Delete()

function Add-AtGoBack([string]$text) {
    $bm = $d.Bookmarks("_GoBack")
    $pt = $d.Range($bm.Start, $bm.Start)
    $pt.InsertBefore($text)
}

# Finish the sentence (merges into the run ending "...make ga") and add
# the closing sentence right after it, still inside the same paragraph.
Add-AtGoBack "meplay more convenient, giving each player their own machine to play on. "
Add-AtGoBack " The rest of this document will outline what such a system will require to be successful.`r"

# Blank line, then the placeholder git-tag marker line.
Add-AtGoBack "`r"
Add-AtGoBack "GIT TAG GOES HERE WHEN COMPLETE --------------------------------------------------------`r"
Add-AtGoBack "`r"
Add-AtGoBack "`r"

# SYSTEM DIAGRAM section.
Add-AtGoBack "SYSTEM DIAGRAM`r"
Add-AtGoBack "Below is a high level diagram of the system:`r"
Add-AtGoBack "`r"
Add-AtGoBack "`r"

# ACTIONS AND THEIR SCENARIOS / PLATFORM / SUMMARY headers.
Add-AtGoBack "ACTIONS AND THEIR SCENARIOS`r"
Add-AtGoBack "`r"
Add-AtGoBack "PLATFORM`r"
Add-AtGoBack "SUMMARY"

# Re-discover the paragraphs we just created (by their text) and apply
# the bold / 18pt (36 half-points) / centered heading formatting used by
# the other section headers in this document.
$headingTexts = @("SYSTEM DIAGRAM", "ACTIONS AND THEIR SCENARIOS", "PLATFORM", "SUMMARY")
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($headingTexts -contains $text) {
        $para.Range.Font.Bold = $true
        $para.Range.Font.Size = 18
        $para.Range.Font.SizeBi = 18
        $para.Format.Alignment = 1
    }
}

# The blank paragraph between "ACTIONS AND THEIR SCENARIOS" and
# "PLATFORM" keeps the same centered / bold / 36 heading formatting too.
$actionsIndex = -1
$platformIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "ACTIONS AND THEIR SCENARIOS") { $actionsIndex = $i }
    if ($text -eq "PLATFORM") { $platformIndex = $i }
}
if ($actionsIndex -gt 0 -and $platformIndex -eq ($actionsIndex + 2)) {
    $blank = $d.Paragraphs.Item($actionsIndex + 1)
    $blank.Range.Font.Bold = $true
    $blank.Range.Font.Size = 18
    $blank.Range.Font.SizeBi = 18
    $blank.Format.Alignment = 1
}

Write-Host "done"
